# Insert a new data row at row 448 of Sheet1, shifting existing rows
# 448-506 down to 449-507, and populate the new row with the latest
# weekly price record for "Zanahoria" (Vega Modelo de Temuco).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift rows 448..506 down by one row.
$ws.Rows.Item(448).Insert()

# Populate the newly inserted row 448 with the new record.
$ws.Cells.Item(448, 1).Value = 10
$ws.Cells.Item(448, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(448, 3).Value = "La Araucanía"
$ws.Cells.Item(448, 4).Value = 45131
$ws.Cells.Item(448, 5).Value = 9
$ws.Cells.Item(448, 6).Value = 100114013
$ws.Cells.Item(448, 7).Value = "Zanahoria"
$ws.Cells.Item(448, 8).Value = "Sin especificar"
$ws.Cells.Item(448, 9).Value = "Primera"
$ws.Cells.Item(448, 10).Value = 250
$ws.Cells.Item(448, 11).Value = 5000
$ws.Cells.Item(448, 12).Value = 5000
$ws.Cells.Item(448, 13).Value = 5000
$ws.Cells.Item(448, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(448, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(448, 16).Value = 200
$ws.Cells.Item(448, 17).Value = 25
$ws.Cells.Item(448, 18).Value = "Hortaliza"
